$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.05755259013108
$ws.Range("D2").Value = 1.060863525316217
$ws.Range("E2").Value = 1.05364348869042
$ws.Range("F2").Value = 1.069306891305747
$ws.Range("I2").Value = 1.051976358614207
$ws.Range("J2").Value = 1.06254841701045
$ws.Range("K2").Value = 1.063589147592431
$ws.Range("L2").Value = 1.056388890393098
$ws.Range("M2").Value = 1.072009749352477

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.058801790855002
$ws.Range("D3").Value = 1.06186652951089
$ws.Range("E3").Value = 1.05472248399805
$ws.Range("F3").Value = 1.07050746540275
$ws.Range("I3").Value = 1.052416036808773
$ws.Range("J3").Value = 1.063448691681091
$ws.Range("K3").Value = 1.064406184900867
$ws.Range("L3").Value = 1.057280290439265
$ws.Range("M3").Value = 1.073025515116331

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.05960960151591
$ws.Range("D4").Value = 1.062515043672137
$ws.Range("E4").Value = 1.05542034430982
$ws.Range("F4").Value = 1.071284212715246
$ws.Range("I4").Value = 1.052698965749616
$ws.Range("J4").Value = 1.064030198626389
$ws.Range("K4").Value = 1.064933753438697
$ws.Range("L4").Value = 1.057856169056118
$ws.Range("M4").Value = 1.07368210433175

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.059949086695822
$ws.Range("D5").Value = 1.062787561587332
$ws.Range("E5").Value = 1.055713649400141
$ws.Range("F5").Value = 1.071610734183192
$ws.Range("I5").Value = 1.052817533617572
$ws.Range("J5").Value = 1.064274419007901
$ws.Range("K5").Value = 1.065155279148802
$ws.Range("L5").Value = 1.05809805050638
$ws.Range("M5").Value = 1.073957973561702

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.060006080884091
$ws.Range("D6").Value = 1.062833311699997
$ws.Range("E6").Value = 1.055762892242888
$ws.Range("F6").Value = 1.071665557274682
$ws.Range("I6").Value = 1.05283741969458
$ws.Range("J6").Value = 1.06431541035825
$ws.Range("K6").Value = 1.065192458845241
$ws.Range("L6").Value = 1.058138650732803
$ws.Range("M6").Value = 1.074004283835239

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.059614138195502
$ws.Range("D7").Value = 1.062518685529499
$ws.Range("E7").Value = 1.055424263760936
$ws.Range("F7").Value = 1.071288575799493
$ws.Range("I7").Value = 1.052700551533396
$ws.Range("J7").Value = 1.064033462872217
$ws.Range("K7").Value = 1.064936714512018
$ws.Range("L7").Value = 1.057859401942218
$ws.Range("M7").Value = 1.073685791139898

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.0579748696789
$ws.Range("D8").Value = 1.061202598985294
$ws.Range("E8").Value = 1.054008207835361
$ws.Range("F8").Value = 1.069712653702174
$ws.Range("I8").Value = 1.052125276030923
$ws.Range("J8").Value = 1.06285288352938
$ws.Range("K8").Value = 1.0638654997242
$ws.Range("L8").Value = 1.056690333990443
$ws.Range("M8").Value = 1.072353173805439

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.055082276726818
$ws.Range("D9").Value = 1.058879617230986
$ws.Range("E9").Value = 1.051510397225224
$ws.Range("F9").Value = 1.066934792788141
$ws.Range("I9").Value = 1.051099483141501
$ws.Range("J9").Value = 1.060764580479794
$ws.Range("K9").Value = 1.061969321715878
$ws.Range("L9").Value = 1.054623194613806
$ws.Range("M9").Value = 1.069999653091789

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.05315102868824
$ws.Range("D10").Value = 1.057328267207135
$ws.Range("E10").Value = 1.049843376418913
$ws.Range("F10").Value = 1.065082169017709
$ws.Range("I10").Value = 1.050407436042756
$ws.Range("J10").Value = 1.05936691138669
$ws.Range("K10").Value = 1.060699354029712
$ws.Range("L10").Value = 1.053240232595084
$ws.Range("M10").Value = 1.068426979996792

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.052314062191992
$ws.Range("D11").Value = 1.056655854459434
$ws.Range("E11").Value = 1.049121083360409
$ws.Range("F11").Value = 1.064279764697687
$ws.Range("I11").Value = 1.050105816672859
$ws.Range("J11").Value = 1.058760383944182
$ws.Range("K11").Value = 1.060148035291061
$ws.Range("L11").Value = 1.052640216845453
$ws.Range("M11").Value = 1.067745102018123

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.05200306351499
$ws.Range("D12").Value = 1.056405988105177
$ws.Range("E12").Value = 1.048852719939245
$ws.Range("F12").Value = 1.063981682561401
$ws.Range("I12").Value = 1.049993486310334
$ws.Range("J12").Value = 1.058534890994113
$ws.Range("K12").Value = 1.059943036539347
$ws.Range("L12").Value = 1.052417164311868
$ws.Range("M12").Value = 1.067491684669193

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.052069778911976
$ws.Range("D13").Value = 1.056459589948282
$ws.Range("E13").Value = 1.048910288080419
$ws.Range("F13").Value = 1.064045623745111
$ws.Range("I13").Value = 1.050017594946792
$ws.Range("J13").Value = 1.058583269158115
$ws.Range("K13").Value = 1.059987019201782
$ws.Range("L13").Value = 1.05246501801484
$ws.Range("M13").Value = 1.067546049809797

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.052288357249788
$ws.Range("D14").Value = 1.056635202527371
$ws.Range("E14").Value = 1.04909890183007
$ws.Range("F14").Value = 1.064255125832042
$ws.Range("I14").Value = 1.050096537445498
$ws.Range("J14").Value = 1.058741748753042
$ws.Range("K14").Value = 1.060131094415708
$ws.Range("L14").Value = 1.052621782936775
$ws.Range("M14").Value = 1.067724157279264

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.052423015628176
$ws.Range("D15").Value = 1.05674338961579
$ws.Range("E15").Value = 1.049215103430048
$ws.Range("F15").Value = 1.064384202475687
$ws.Range("I15").Value = 1.050145137329851
$ws.Range("J15").Value = 1.058839366474839
$ws.Range("K15").Value = 1.060219835438748
$ws.Range("L15").Value = 1.052718347069311
$ws.Range("M15").Value = 1.067833876880446

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.053206559869739
$ws.Range("D16").Value = 1.057372878785884
$ws.Range("E16").Value = 1.049891302707825
$ws.Range("F16").Value = 1.06513541739677
$ws.Range("I16").Value = 1.050427412157219
$ws.Range("J16").Value = 1.059407136458083
$ws.Range("K16").Value = 1.060735913274971
$ws.Range("L16").Value = 1.053280028540303
$ws.Range("M16").Value = 1.068472214833851

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.053697860426765
$ws.Range("D17").Value = 1.057767560494539
$ws.Range("E17").Value = 1.050315339393955
$ws.Range("F17").Value = 1.065606577757962
$ws.Range("I17").Value = 1.050603950481245
$ws.Range("J17").Value = 1.059762926599857
$ws.Range("K17").Value = 1.061059255241269
$ws.Range("L17").Value = 1.053632038014997
$ws.Range("M17").Value = 1.068872384771209

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.05398435816645
$ws.Range("D18").Value = 1.057997707344724
$ws.Range("E18").Value = 1.050562628388061
$ws.Range("F18").Value = 1.065881377919236
$ws.Range("I18").Value = 1.050706733427319
$ws.Range("J18").Value = 1.059970325035992
$ws.Range("K18").Value = 1.06124771880737
$ws.Range("L18").Value = 1.053837245082009
$ws.Range("M18").Value = 1.069105710428633

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.054082034794439
$ws.Range("D19").Value = 1.058076170623506
$ws.Range("E19").Value = 1.050646940074339
$ws.Range("F19").Value = 1.065975074405154
$ws.Range("I19").Value = 1.050741747758021
$ws.Range("J19").Value = 1.060041020909583
$ws.Range("K19").Value = 1.061311956946803
$ws.Range("L19").Value = 1.053907196097112
$ws.Range("M19").Value = 1.069185253810991

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.053645155747604
$ws.Range("D20").Value = 1.057725221549867
$ws.Range("E20").Value = 1.050269848883252
$ws.Range("F20").Value = 1.065556028767724
$ws.Range("I20").Value = 1.050585029135469
$ws.Range("J20").Value = 1.059724766934606
$ws.Range("K20").Value = 1.061024577824757
$ws.Range("L20").Value = 1.053594282548039
$ws.Range("M20").Value = 1.068829459309018

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.052223994497659
$ws.Range("D21").Value = 1.05658349184026
$ws.Range("E21").Value = 1.049043361776632
$ws.Range("F21").Value = 1.064193433632469
$ws.Range("I21").Value = 1.050073299014655
$ws.Range("J21").Value = 1.058695086046713
$ws.Range("K21").Value = 1.06008867379093
$ws.Range("L21").Value = 1.05257562455924
$ws.Range("M21").Value = 1.067671712882853

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.051329804021972
$ws.Range("D22").Value = 1.055865048080421
$ws.Range("E22").Value = 1.048271805191484
$ws.Range("F22").Value = 1.063336521246399
$ws.Range("I22").Value = 1.049749843679928
$ws.Range("J22").Value = 1.058046517006922
$ws.Range("K22").Value = 1.059498992220821
$ws.Range("L22").Value = 1.051934111576795
$ws.Range("M22").Value = 1.066942995051001

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.051803894054776
$ws.Range("D23").Value = 1.056245965578082
$ws.Range("E23").Value = 1.048680862054106
$ws.Range("F23").Value = 1.06379080601091
$ws.Range("I23").Value = 1.049921476014494
$ws.Range("J23").Value = 1.058390447160202
$ws.Range("K23").Value = 1.059811711911691
$ws.Range("L23").Value = 1.052274289237205
$ws.Range("M23").Value = 1.06732937843995

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.053668970938335
$ws.Range("D24").Value = 1.057744352895255
$ws.Range("E24").Value = 1.050290404224228
$ws.Range("F24").Value = 1.065578869741459
$ws.Range("I24").Value = 1.050593579460739
$ws.Range("J24").Value = 1.059742010040051
$ws.Range("K24").Value = 1.061040247477851
$ws.Range("L24").Value = 1.053611342970953
$ws.Range("M24").Value = 1.068848855745558

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.055830571778574
$ws.Range("D25").Value = 1.059480631872241
$ws.Range("E25").Value = 1.052156452988993
$ws.Range("F25").Value = 1.067653053767673
$ws.Range("I25").Value = 1.051366113430412
$ws.Range("J25").Value = 1.061305412576682
$ws.Range("K25").Value = 1.062460553251611
$ws.Range("L25").Value = 1.055158450506855
$ws.Range("M25").Value = 1.070608731483715
